# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H1 header: copy the header formatting from the existing "sum" header (G1)
# so it picks up the same bold/bordered/centered style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# H2 / H3 data cells: numeric 0, default (unstyled) like the other data cells.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
